$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-11 held "Surfboards" (rows 2-5) and "Surfboard Accessories" (rows 6-11)
# in the Product Group column (D). Both labels are consolidated into a single
# "Surfing" category.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = "Surfing"
}

# Restore the view: scrolled to the top (A2 as the frozen pane's top-left
# cell) with D12 selected.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()
